# Update Name of Algo
# Apply updated KNN imputation result values to the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = -7.907999999999999
$ws.Range("B7").Value = 5.941000000000001
$ws.Range("A9").Value = -21.752
$ws.Range("B12").Value = 5.927
$ws.Range("B14").Value = 6.225
$ws.Range("D15").Value = -8.231999999999999
$ws.Range("A18").Value = -22.098
$ws.Range("A20").Value = -20.457
$ws.Range("B26").Value = 5.761
$ws.Range("A27").Value = -21.188
$ws.Range("B27").Value = 5.695000000000001
$ws.Range("B29").Value = 5.958
$ws.Range("D33").Value = -7.845000000000001
$ws.Range("A35").Value = -19.873
$ws.Range("D35").Value = -7.551
$ws.Range("B37").Value = 8.847000000000001
$ws.Range("B38").Value = 6.273
$ws.Range("D38").Value = -8.669
$ws.Range("D43").Value = -7.870000000000002
$ws.Range("D44").Value = -7.741
$ws.Range("D47").Value = -7.848999999999999
$ws.Range("B51").Value = 5.465
$ws.Range("D51").Value = -8.190000000000001
$ws.Range("B52").Value = 5.337000000000001
$ws.Range("B55").Value = 5.705
$ws.Range("D57").Value = -8.038
$ws.Range("D63").Value = -7.336999999999999
$ws.Range("A69").Value = -21.586
$ws.Range("B69").Value = 5.952999999999999
$ws.Range("B70").Value = 5.412
$ws.Range("D70").Value = -6.797
$ws.Range("A76").Value = -20.66
$ws.Range("A78").Value = -19.993
$ws.Range("B81").Value = 6.02
$ws.Range("A82").Value = -22.156
$ws.Range("A83").Value = -20.469
$ws.Range("B83").Value = 7.776999999999999
$ws.Range("D88").Value = -8.055
$ws.Range("A93").Value = -21.958
$ws.Range("D99").Value = -8.217000000000001
$ws.Range("B102").Value = 7.575
